{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the content of: async (context) => { ... }\n//\n// Change description: the placeholder line\n//   \">>>  your stuff after this line >>>\"\n// is replaced with\n//   \">>>  MY CHANGES FOR ASSIGNMENT 1 >>>\"\n// (the surrounding \">>>  \" ... \" >>>\" wrapper is kept unchanged).\n\nconst body = context.document.body;\n\n// Search for the exact phrase that needs to be swapped out. Using the\n// Word.js search API keeps this robust even though the phrase happens to\n// span several runs (and a couple of proofing-error marks) in the\n// original document.\nconst results = body.search(\"your stuff after this line\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n    results.items[0].insertText(\"MY CHANGES FOR ASSIGNMENT 1\", Word.InsertLocation.replace);\n    await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change description: the placeholder line\n#   \">>>  your stuff after this line >>>\"\n# is replaced with\n#   \">>>  MY CHANGES FOR ASSIGNMENT 1 >>>\"\n# (the surrounding \">>>  \" ... \" >>>\" wrapper is kept unchanged).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# wdFindContinue = 1, wdReplaceOne = 1 (explicit arg below), wdReplaceAll = 2\n$find.Execute(\n    \"your stuff after this line\",  # FindText\n    $false,                         # MatchCase\n    $false,                         # MatchWholeWord\n    $false,                         # MatchWildcards\n    $false,                         # MatchSoundsLike\n    $false,                         # MatchAllWordForms\n    $true,                          # Forward\n    1,                              # Wrap (wdFindContinue)\n    $false,                         # Format\n    \"MY CHANGES FOR ASSIGNMENT 1\",  # ReplaceWith\n    2                               # Replace (wdReplaceAll)\n)\n"}
